# Auto-generated edit script: updates crypto price/volume table to refreshed values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.539.06"
$ws.Range("E2").Value = "  -1.27%  "
$ws.Range("D3").Value = "2.220.07"
$ws.Range("E3").Value = "  -2.21%  "
$ws.Range("E4").Value = "  +0.29%  "
$ws.Range("D5").Value = "'312.34"
$ws.Range("E5").Value = "  -2.10%  "
$ws.Range("D6").Value = "'97.40"
$ws.Range("E6").Value = "  -4.96%  "
$ws.Range("D7").Value = "'0.567"
$ws.Range("E7").Value = "  -3.30%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("D9").Value = "'0.533"
$ws.Range("E9").Value = "  -6.74%  "
$ws.Range("D10").Value = "'35.69"
$ws.Range("E10").Value = "  -7.84%  "
$ws.Range("D11").Value = "'0.0821"
$ws.Range("E11").Value = "  -2.16%  "
$ws.Range("D12").Value = "'7.35"
$ws.Range("E12").Value = "  -6.54%  "
$ws.Range("E13").Value = "  -3.32%  "
$ws.Range("D14").Value = "2.558.86"
$ws.Range("E14").Value = "  -1.97%  "
$ws.Range("D15").Value = "2.228.63"
$ws.Range("E15").Value = "  -1.62%  "
$ws.Range("D16").Value = "'0.836"
$ws.Range("E16").Value = "  -4.53%  "
$ws.Range("D17").Value = "'14.05"
$ws.Range("E17").Value = "  -3.57%  "
$ws.Range("D18").Value = "43.456.80"
$ws.Range("E18").Value = "  -1.31%  "
$ws.Range("D19").Value = "'12.94"
$ws.Range("E19").Value = "  -10.55%  "
$ws.Range("D20").Value = "0.0₃0964"
$ws.Range("E20").Value = "  -3.02%  "
$ws.Range("D21").Value = "'6.27"
$ws.Range("E21").Value = "  -5.97%  "
$ws.Range("D22").Value = "'65.16"
$ws.Range("E22").Value = "  -1.42%  "
$ws.Range("D23").Value = "'234.42"
$ws.Range("E23").Value = "  -2.00%  "
$ws.Range("D24").Value = "'2.96"
$ws.Range("E24").Value = "  -7.64%  "
$ws.Range("D25").Value = "'2.02"
$ws.Range("E25").Value = "  -7.89%  "
$ws.Range("E26").Value = "  +0.19%  "
$ws.Range("D27").Value = "'9.98"
$ws.Range("E27").Value = "  -2.66%  "
$ws.Range("D28").Value = "'2.21"
$ws.Range("E28").Value = "  -0.46%  "
$ws.Range("D29").Value = "'36.16"
$ws.Range("E29").Value = "  -7.76%  "
$ws.Range("D30").Value = "'160.73"
$ws.Range("E30").Value = "  -1.86%  "
$ws.Range("D31").Value = "'5.92"
$ws.Range("E31").Value = "  -9.15%  "
$ws.Range("D32").Value = "'19.81"
$ws.Range("E32").Value = "  -3.24%  "
$ws.Range("D33").Value = "'0.0824"
$ws.Range("E33").Value = "  -6.61%  "
$ws.Range("E34").Value = "  -1.20%  "
$ws.Range("D35").Value = "'3.12"
$ws.Range("E35").Value = "  -4.04%  "
$ws.Range("E36").Value = "  +0.59%  "
$ws.Range("D37").Value = "'1.86"
$ws.Range("E37").Value = "  -8.53%  "
$ws.Range("D38").Value = "'0.116"
$ws.Range("E38").Value = "  -3.72%  "
$ws.Range("D39").Value = "'15.34"
$ws.Range("E39").Value = "  -3.07%  "
$ws.Range("D40").Value = "'3.53"
$ws.Range("E40").Value = "  -8.82%  "
$ws.Range("D41").Value = "'3.96"
$ws.Range("E41").Value = "  -13.37%  "
$ws.Range("D42").Value = "'0.0306"
$ws.Range("E42").Value = "  -6.24%  "
$ws.Range("D43").Value = "'1.01"
$ws.Range("E43").Value = "  +0.36%  "
$ws.Range("D44").Value = "1.702.25"
$ws.Range("E44").Value = "  -4.19%  "
$ws.Range("D45").Value = "'82.56"
$ws.Range("E45").Value = "  -3.17%  "
$ws.Range("D46").Value = "'0.193"
$ws.Range("E46").Value = "  -7.15%  "
$ws.Range("D47").Value = "'5.09"
$ws.Range("E47").Value = "  -5.96%  "
$ws.Range("D48").Value = "'71.83"
$ws.Range("E48").Value = "  -4.28%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Value = "'101.23"
$ws.Range("E49").Value = "  -3.54%  "
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").Value = "'1.64"
$ws.Range("E50").Value = "  +0.53%  "
$ws.Range("D51").Value = "'56.25"
$ws.Range("E51").Value = "  -5.93%  "
